$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.221213
$ws.Range("N2").Value = 0.6636390000000001
$ws.Range("O2").Value = 0.1279649941053948
$ws.Range("P2").Value = 0.1279649941053949
$ws.Range("Q2").Value = 0.184275959325
$ws.Range("R2").Value = 1.658483633925
$ws.Range("S2").Value = 0.1279649941053948
$ws.Range("T2").Value = 0.1279649941053949

# Row 3 updates
$ws.Range("O3").Value = 0.6284678770050237
$ws.Range("P3").Value = 0.6284678770050239
$ws.Range("S3").Value = 0.6284678770050237
$ws.Range("T3").Value = 0.6284678770050239

# Row 4 updates
$ws.Range("O4").Value = 0.2435671288895813
$ws.Range("P4").Value = 0.2435671288895814
$ws.Range("S4").Value = 0.2435671288895813
$ws.Range("T4").Value = 0.2435671288895814
